$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Insert a new "Texas Notes" worksheet between "Data" and "PPEIdtICEaT"
# ------------------------------------------------------------------
$dataSheet = $wb.Worksheets.Item("Data")
$notesSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$notesSheet.Name = "Texas Notes"

# Populate it with the reviewer's notes (column A, rows 1,3,4,5,7,9,10,12)
$notesSheet.Range("A1").Value = "This spreadsheet uses a very particular study. "
$notesSheet.Range("A3").Value = "It's done in Wisconsin where the authors use an educational seminar for builders"
$notesSheet.Range("A4").Value = "then they follow up with phone surveys to see what the builders actually implemented"
$notesSheet.Range("A5").Value = "then they use those results to try and estimate how much energy those builders decisions saved"
$notesSheet.Range("A7").Value = "It's all very niche and a bit subjective"
$notesSheet.Range("A9").Value = "That said, I did a quick literature search and didn't find anything that I thought"
$notesSheet.Range("A10").Value = "would give us better or more Texas-specific numbers. "
$notesSheet.Range("A12").Value = "Since this is a `"low`" priority sheet, I will leave it alone."

# ------------------------------------------------------------------
# Restore each sheet's remembered selection (and, by selecting last,
# make "PPEIdtICEaT" the active / tab-selected sheet again)
# ------------------------------------------------------------------
$aboutSheet = $wb.Worksheets.Item("About")
$aboutSheet.Range("B25").Select() | Out-Null

$dataSheet.Range("A14").Select() | Out-Null

$notesSheet.Range("C27").Select() | Out-Null

$ppeSheet = $wb.Worksheets.Item("PPEIdtICEaT")
$ppeSheet.Range("D20").Select() | Out-Null
